{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of: async (context) => { ... }\n\nconst newTexts = [\n  \"1) Analyze the Pestel Open AI\",\n  \"Introduction\",\n  \"Recently, Open AI has unveiled a new system based on artificial intelligence, Chat GPT 4.O; Artificial Intelligence has developed greatly in recent years. In 2023, the company employs 375 people, with a turnover exceeding $2 billion, and this is thanks to very strong growth year after year. The company specializes in research and development of artificial intelligence worldwide. The main objective is to be able to benefit everyone from the benefits of this technology while gradually reducing the constraints and/or threats associated with it. The application allows among other things to generate text and images for individuals, but also on behalf of companies.\",\n  \"In this article, we will study the PESTEL analysis of Open AI. This analysis highlights the main environmental factors related to the company in the political, economic, sociological, technological, ecological and legal fields.\",\n  \"1) PESTEL OPEN AI is a type of artificial intelligence that uses the PESTEL framework to analyze and predict trends in business, economics, politics, technology, environment, and society. It can be used by businesses to make informed decisions about their operations and strategies.\",\n  \"Politics\",\n  \"There are many political factors that come into play in the development of this company. Factors that are governmental first, since the various governments in place have a considerable influence on research and development of artificial intelligence. Moreover, international relations must be conducive to understanding because Open AI wants to be able to surround itself with researchers from around the world in order to expand its range of skills. If there are commercial disputes between countries, this necessarily affects the company. Open AI is developing thanks to private and government funding.\",\n  \"1) economic\",\n  \"Open AI is also very impacted by various economic factors. First of all, investments at this level must be very important because research can be very expensive in the field of artificial intelligence. In addition, economic growth is also to be taken into account, trends, people who tend to be more and more interested in new technologies and everything they bring them on a daily basis. More people who are interested in the sector, it's also more chance to recruit.\",\n  \"In 2024, the main competitors of Open AI are Whatsthebigdata, DataCamp, Writesonic and Anakin.ai.\",\n  \"Sociological\",\n  \"'OpenAI is a group that wants to develop products that can be useful to everyone. The demand changes depending on the population, for example, AI is used to develop care for the elderly, to help children and parents at the educational level. Research is necessary in order to have better adaptation. There are also cultural values \\u200b\\u200band the notion of ethics to take into account. Consumer attitudes are changing, mobile applications are needed for OpenAI because more and more people are on their mobile all day or almost.'\",\n  \"Technological\",\n  \"Of course, Open AI is directly linked to new technologies, and research and development are more precisely related to computing power, necessary in order to be able to process as much data as possible in a minimum of time.\",\n  \"Cloud computing can significantly improve the profitability of products and business growth drivers. Among these products, ChatGPT is available for both individuals and businesses. GPT 4-0 is a new service that has been active since May 13th, 2024. There are also other options such as GPT 4 Turbo or GPT 3.5 Turbo. These applications allow you to generate text and/or images using resources found on the internet.\",\n  \"'Ecological' is the correct translation of '\u00c9cologique'.\",\n  \"Such applications are very energy hungry, so they are harmful to the environment. It is important for Open AI to find more efficient ways to reduce these energy costs. Among these means, the possibility of finding partnerships with eco-responsible companies, having effective waste management, or even using recycling operations in order to always be more in line with environmental standards.\",\n  \"Legal\",\n  \"At the legal level, Open Ai is subject to a large number of regulations. The company must file patents in case of innovations, especially given the extent of competition, which has become more and more important over the years. There is also the notion of intellectual property and data confidentiality that comes into play. Open AI is subject to the General Data Protection Regulation (GDPR) at the European level as well as the California Consumer Privacy Act at the US level.\",\n  \"1) CONCLUSION = CONCLUSION\",\n  \"OpenAI is a company that has come a long way since its inception in 2015. Today, it wants to expand across the planet to serve everyone's interests and reduce anything that could lead to poor management and use of AI. Research is plentiful, and the group's revenue is constantly increasing.\",\n  \"\u00c9cologique\",\n  \"De telles applications sont tr\u00e8s gourmandes en \u00e9nergie, elles sont donc n\u00e9fastes pour l\\u2019environnement. Il est important pour Open AI de trouver des moyens plus efficaces pour diminuer ces d\u00e9penses \u00e9nerg\u00e9tiques. Parmi ces moyens, la possibilit\u00e9 de trouver des partenariats avec des entreprises \u00e9co-responsables, avoir une gestion des d\u00e9chets efficace, ou encore utiliser des op\u00e9rations de recyclage afin d\\u2019\u00eatre toujours plus en phase avec les normes environnementales.\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst keepCount = newTexts.length;\n\n// Replace the text of the paragraphs we keep (indices 0..22).\nfor (let i = 0; i < keepCount; i++) {\n  paragraphs.items[i].insertText(newTexts[i], \"Replace\");\n}\nawait context.sync();\n\n// Delete the trailing paragraphs that no longer exist in the target\n// document (original indices 23..30). Delete from the end backwards so\n// earlier indices stay valid while deleting.\nfor (let i = paragraphs.items.length - 1; i >= keepCount; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"1) Analyze the Pestel Open AI\",\n  \"Introduction\",\n  \"Recently, Open AI has unveiled a new system based on artificial intelligence, Chat GPT 4.O; Artificial Intelligence has developed greatly in recent years. In 2023, the company employs 375 people, with a turnover exceeding `$2 billion, and this is thanks to very strong growth year after year. The company specializes in research and development of artificial intelligence worldwide. The main objective is to be able to benefit everyone from the benefits of this technology while gradually reducing the constraints and/or threats associated with it. The application allows among other things to generate text and images for individuals, but also on behalf of companies.\",\n  \"In this article, we will study the PESTEL analysis of Open AI. This analysis highlights the main environmental factors related to the company in the political, economic, sociological, technological, ecological and legal fields.\",\n  \"1) PESTEL OPEN AI is a type of artificial intelligence that uses the PESTEL framework to analyze and predict trends in business, economics, politics, technology, environment, and society. It can be used by businesses to make informed decisions about their operations and strategies.\",\n  \"Politics\",\n  \"There are many political factors that come into play in the development of this company. Factors that are governmental first, since the various governments in place have a considerable influence on research and development of artificial intelligence. Moreover, international relations must be conducive to understanding because Open AI wants to be able to surround itself with researchers from around the world in order to expand its range of skills. If there are commercial disputes between countries, this necessarily affects the company. Open AI is developing thanks to private and government funding.\",\n  \"1) economic\",\n  \"Open AI is also very impacted by various economic factors. First of all, investments at this level must be very important because research can be very expensive in the field of artificial intelligence. In addition, economic growth is also to be taken into account, trends, people who tend to be more and more interested in new technologies and everything they bring them on a daily basis. More people who are interested in the sector, it's also more chance to recruit.\",\n  \"In 2024, the main competitors of Open AI are Whatsthebigdata, DataCamp, Writesonic and Anakin.ai.\",\n  \"Sociological\",\n  \"'OpenAI is a group that wants to develop products that can be useful to everyone. The demand changes depending on the population, for example, AI is used to develop care for the elderly, to help children and parents at the educational level. Research is necessary in order to have better adaptation. There are also cultural values \u200b\u200band the notion of ethics to take into account. Consumer attitudes are changing, mobile applications are needed for OpenAI because more and more people are on their mobile all day or almost.'\",\n  \"Technological\",\n  \"Of course, Open AI is directly linked to new technologies, and research and development are more precisely related to computing power, necessary in order to be able to process as much data as possible in a minimum of time.\",\n  \"Cloud computing can significantly improve the profitability of products and business growth drivers. Among these products, ChatGPT is available for both individuals and businesses. GPT 4-0 is a new service that has been active since May 13th, 2024. There are also other options such as GPT 4 Turbo or GPT 3.5 Turbo. These applications allow you to generate text and/or images using resources found on the internet.\",\n  \"'Ecological' is the correct translation of '\u00c9cologique'.\",\n  \"Such applications are very energy hungry, so they are harmful to the environment. It is important for Open AI to find more efficient ways to reduce these energy costs. Among these means, the possibility of finding partnerships with eco-responsible companies, having effective waste management, or even using recycling operations in order to always be more in line with environmental standards.\",\n  \"Legal\",\n  \"At the legal level, Open Ai is subject to a large number of regulations. The company must file patents in case of innovations, especially given the extent of competition, which has become more and more important over the years. There is also the notion of intellectual property and data confidentiality that comes into play. Open AI is subject to the General Data Protection Regulation (GDPR) at the European level as well as the California Consumer Privacy Act at the US level.\",\n  \"1) CONCLUSION = CONCLUSION\",\n  \"OpenAI is a company that has come a long way since its inception in 2015. Today, it wants to expand across the planet to serve everyone's interests and reduce anything that could lead to poor management and use of AI. Research is plentiful, and the group's revenue is constantly increasing.\",\n  \"\u00c9cologique\",\n  \"De telles applications sont tr\u00e8s gourmandes en \u00e9nergie, elles sont donc n\u00e9fastes pour l\u2019environnement. Il est important pour Open AI de trouver des moyens plus efficaces pour diminuer ces d\u00e9penses \u00e9nerg\u00e9tiques. Parmi ces moyens, la possibilit\u00e9 de trouver des partenariats avec des entreprises \u00e9co-responsables, avoir une gestion des d\u00e9chets efficace, ou encore utiliser des op\u00e9rations de recyclage afin d\u2019\u00eatre toujours plus en phase avec les normes environnementales.\"\n)\n\n$keepCount = $newTexts.Count\n\n# Replace the text of the paragraphs we keep (1-based indices 1..23).\nfor ($i = 0; $i -lt $keepCount; $i++) {\n    $d.Paragraphs($i + 1).Range.Text = $newTexts[$i]\n}\n\n# Delete the trailing paragraphs that no longer exist in the target\n# document (original 1-based indices 24..31). Walk from the end backwards\n# so earlier indices stay valid while deleting.\nfor ($i = $d.Paragraphs.Count; $i -gt $keepCount; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
